$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$row = 49

# Columns A-L hold text values (some look numeric/date-like, e.g. "-497",
# "12", "0", "7/2/2025"); force them to be stored as text so Excel does not
# auto-convert them into numbers or dates, matching the source data which
# used inline strings for every cell in A:L.
$textRange = $ws.Range("A" + $row + ":L" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "-497"
$ws.Cells.Item($row, 2).Value = "7/2/2025"
$ws.Cells.Item($row, 3).Value = "Machain 4556"
$ws.Cells.Item($row, 4).Value = "12"
$ws.Cells.Item($row, 5).Value = "807918340"
$ws.Cells.Item($row, 6).Value = "NEW"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Poste"
$ws.Cells.Item($row, 9).Value = "0"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Poste"

# Restore the default (unstyled) cell style so the new cells match the rest
# of the sheet, which carries no explicit style index.
$textRange.Style = "Normal"

# Columns M and N hold numeric coordinates.
$ws.Cells.Item($row, 13).Value = -58.492573
$ws.Cells.Item($row, 14).Value = -34.551355
